$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Trim trailing whitespace from a handful of shared strings ---
$ws.Range("D2").Value = "Amount"
$ws.Range("C11").Value = "1 MOhm"
$ws.Range("C12").Value = "3 MOhm"
$ws.Range("H20:H25").Value = "Kan goedkoper maar stekkerbaar is handig"

# --- Battery-management-disable rows: drop the "Amount" entirely (was 1,1,5,..,2) ---
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D26").ClearContents()

# --- Other rows: amount explicitly set to 0 instead of removed ---
$ws.Range("D8").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("D17").Value = 0

# --- New blank cell at C1 (keeps dimension/formatting in step with the rest of row 1) ---
$ws.Range("C1").Value = ""

# --- New "current gear" estimation rows 37-63, Amount-only column D ---
$gearAmounts = @(1,1,1,5,1,5,1,1,1,1,2,1,1,1,1,1,1,1,1,1,1,2,1,1,2)
$r = 37
foreach ($amt in $gearAmounts) {
    $ws.Cells.Item($r, 4).Value = $amt
    $r++
}

# --- Move the active selection to reflect where the author was last working ---
$ws.Range("H40").Select()

Write-Output "edit applied"
